$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("customer-rules")

# JIRA:101 - defect fix: make the ACTION code null-safe so that when the
# output already has content, new action text is appended with a comma
# separator instead of being blindly concatenated (which could NPE /
# produce a garbled string when orderOutput.getOutput() is null or empty).
$ws.Range("E8").Value = "System.out.println(`"Output : `" + `$1 );`norderOutput.setOutput( (orderOutput.getOutput()  != null && !orderOutput.getOutput().trim() .equals(`"`") )  ? orderOutput.getOutput() +  `" ,`" +`$1 : `$1);"

# The sheet was left scrolled/selected further down while editing (D3/E11);
# restore the view to the top of the table with E8 (the action cell that
# was just changed) selected.
$ws.Range("E8").Select()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 4
